$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

$ws.Range("C53").Value = "so far we only use the full forms (can, may, might, could), while Biber also includes the contractions"
$ws.Range("D53").Value = "works well"

$ws.Range("C54").Value = "so far we only use the full forms (ought, should, must), while Biber also includes the contractions"
$ws.Range("D54").Value = "works well"

$ws.Range("D55").Value = "works well"

$ws.Range("D60").Value = "works well, but open question whether tagger will remove apostrophes as the beginning of clitics? (we currently rely on that for identificaiton)"

$ws.Range("C64").Value = "Biber only allows for one or two intervening adverbs, we allow  for more"
$ws.Range("D64").Value = "works well"

$ws.Activate()
$ws.Range("C65").Select()
$excel.ActiveWindow.ScrollRow = 63
